# Daily attendance processing - 2025-10-18 02:42:33
# Normalize the "Recorded By" (column G) values: when the list of recorders
# starts with "System", move it so the swap is reflected by exchanging the
# first and last comma-separated entries (e.g. "System, a, b" -> "b, a, System",
# and "System, a" -> "a, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val.GetType().Name -eq "String" -and $val.StartsWith("System, ")) {
        $parts = $val -split ", "
        if ($parts.Count -ge 2) {
            $tmp = $parts[0]
            $parts[0] = $parts[$parts.Count - 1]
            $parts[$parts.Count - 1] = $tmp
            $cell.Value = $parts -join ", "
        }
    }
}
